# Natmi following Dr Hou advice
# Rebuild the L1cam-Erbb3 ligand-receptor table: expand from 8 data rows (2 sending
# clusters x 4 combos) to 16 data rows (4 sending clusters x 4 target clusters),
# with refreshed expression/specificity statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (L1cam/Erbb3)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "L1cam"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04178033333333334
$ws.Range("N2").Value = 0.125341
$ws.Range("O2").Value = 0.009630623139527362
$ws.Range("P2").Value = 0.009630623139527362
$ws.Range("Q2").Value = 0.753941977599889
$ws.Range("R2").Value = 6.785477798399001
$ws.Range("S2").Value = 0.006546858830902224
$ws.Range("T2").Value = 0.006546858830902224

# Row 3: ECs -> FAPs (L1cam/Erbb3)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "L1cam"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.181585
$ws.Range("N3").Value = 0.544755
$ws.Range("O3").Value = 0.04185645645377991
$ws.Range("P3").Value = 0.04185645645377991
$ws.Range("Q3").Value = 3.276770266771667
$ws.Range("R3").Value = 29.490932400945
$ws.Range("S3").Value = 0.02845385055511078
$ws.Range("T3").Value = 0.02845385055511078

# Row 4: ECs -> M2 (L1cam/Erbb3)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "L1cam"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4121513333333333
$ws.Range("N4").Value = 1.236454
$ws.Range("O4").Value = 0.0950034107224385
$ws.Range("P4").Value = 0.0950034107224385
$ws.Range("Q4").Value = 7.43742729012289
$ws.Range("R4").Value = 66.93684561110601
$ws.Range("S4").Value = 0.06458293606165881
$ws.Range("T4").Value = 0.06458293606165881

# Row 5: ECs -> sCs (L1cam/Erbb3)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "L1cam"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.04537966666667
$ws.Range("H5").Value = 54.13613900000001
$ws.Range("I5").Value = 0.6797959733292525
$ws.Range("J5").Value = 0.6797959733292525
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.702762666666667
$ws.Range("N5").Value = 11.108288
$ws.Range("O5").Value = 0.8535095096842542
$ws.Range("P5").Value = 0.8535095096842543
$ws.Range("Q5").Value = 66.81775813555912
$ws.Range("R5").Value = 601.3598232200321
$ws.Range("S5").Value = 0.5802123278815806
$ws.Range("T5").Value = 0.5802123278815807

# Row 6: FAPs -> ECs (L1cam/Erbb3)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "L1cam"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04178033333333334
$ws.Range("N6").Value = 0.125341
$ws.Range("O6").Value = 0.009630623139527362
$ws.Range("P6").Value = 0.009630623139527362
$ws.Range("Q6").Value = 0.02507652821311112
$ws.Range("R6").Value = 0.225688753918
$ws.Range("S6").Value = 0.0002177521547520203
$ws.Range("T6").Value = 0.0002177521547520203

# Row 7: FAPs -> FAPs (L1cam/Erbb3)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "L1cam"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.181585
$ws.Range("N7").Value = 0.544755
$ws.Range("O7").Value = 0.04185645645377991
$ws.Range("P7").Value = 0.04185645645377991
$ws.Range("Q7").Value = 0.1089871959433334
$ws.Range("R7").Value = 0.9808847634900001
$ws.Range("S7").Value = 0.0009463908462668786
$ws.Range("T7").Value = 0.0009463908462668786

# Row 8: FAPs -> M2 (L1cam/Erbb3)
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "L1cam"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.02261039099934159
$ws.Range("J8").Value = 0.02261039099934159
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4121513333333333
$ws.Range("N8").Value = 1.236454
$ws.Range("O8").Value = 0.0950034107224385
$ws.Range("P8").Value = 0.0950034107224385
$ws.Range("Q8").Value = 0.2473729554991111
$ws.Range("R8").Value = 2.226356599492
$ws.Range("S8").Value = 0.002148064262705376
$ws.Range("T8").Value = 0.002148064262705376

# Row 9: FAPs -> sCs (L1cam/Erbb3)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "L1cam"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.02261039099934159
$ws.Range("J9").Value = 0.02261039099934159
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.702762666666667
$ws.Range("N9").Value = 11.108288
$ws.Range("O9").Value = 0.8535095096842542
$ws.Range("P9").Value = 0.8535095096842543
$ws.Range("Q9").Value = 2.222395684024889
$ws.Range("R9").Value = 20.001561156224
$ws.Range("S9").Value = 0.01929818373561732
$ws.Range("T9").Value = 0.01929818373561732

# Row 10: M2 -> ECs (L1cam/Erbb3)
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "L1cam"
$ws.Range("C10").Value = "Erbb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.04178033333333334
$ws.Range("N10").Value = 0.125341
$ws.Range("O10").Value = 0.009630623139527362
$ws.Range("P10").Value = 0.009630623139527362
$ws.Range("Q10").Value = 0.1959520333981111
$ws.Range("R10").Value = 1.763568300583
$ws.Range("S10").Value = 0.001701550435445419
$ws.Range("T10").Value = 0.001701550435445419

# Row 11: M2 -> FAPs (L1cam/Erbb3)
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "L1cam"
$ws.Range("C11").Value = "Erbb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.690054333333333
$ws.Range("H11").Value = 14.070163
$ws.Range("I11").Value = 0.1766812397072912
$ws.Range("J11").Value = 0.1766812397072912
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.181585
$ws.Range("N11").Value = 0.544755
$ws.Range("O11").Value = 0.04185645645377991
$ws.Range("P11").Value = 0.04185645645377991
$ws.Range("Q11").Value = 0.8516435161183333
$ws.Range("R11").Value = 7.664791645065001
$ws.Range("S11").Value = 0.007395250616008082
$ws.Range("T11").Value = 0.007395250616008084

# Row 12: M2 -> M2 (L1cam/Erbb3)
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "L1cam"
$ws.Range("C12").Value = "Erbb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.690054333333333
$ws.Range("H12").Value = 14.070163
$ws.Range("I12").Value = 0.1766812397072912
$ws.Range("J12").Value = 0.1766812397072912
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.4121513333333333
$ws.Range("N12").Value = 1.236454
$ws.Range("O12").Value = 0.0950034107224385
$ws.Range("P12").Value = 0.0950034107224385
$ws.Range("Q12").Value = 1.933012146889111
$ws.Range("R12").Value = 17.397109322002
$ws.Range("S12").Value = 0.01678532038286139
$ws.Range("T12").Value = 0.01678532038286139

# Row 13: M2 -> sCs (L1cam/Erbb3)
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "L1cam"
$ws.Range("C13").Value = "Erbb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.690054333333333
$ws.Range("H13").Value = 14.070163
$ws.Range("I13").Value = 0.1766812397072912
$ws.Range("J13").Value = 0.1766812397072912
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.702762666666667
$ws.Range("N13").Value = 11.108288
$ws.Range("O13").Value = 0.8535095096842542
$ws.Range("P13").Value = 0.8535095096842543
$ws.Range("Q13").Value = 17.36615809010489
$ws.Range("R13").Value = 156.295422810944
$ws.Range("S13").Value = 0.1507991182729763
$ws.Range("T13").Value = 0.1507991182729763

# Row 14: sCs -> ECs (L1cam/Erbb3)
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "L1cam"
$ws.Range("C14").Value = "Erbb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.209654333333333
$ws.Range("H14").Value = 9.628962999999999
$ws.Range("I14").Value = 0.1209123959641148
$ws.Range("J14").Value = 0.1209123959641148
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.04178033333333334
$ws.Range("N14").Value = 0.125341
$ws.Range("O14").Value = 0.009630623139527362
$ws.Range("P14").Value = 0.009630623139527362
$ws.Range("Q14").Value = 0.1341004279314444
$ws.Range("R14").Value = 1.206903851383
$ws.Range("S14").Value = 0.001164461718427699
$ws.Range("T14").Value = 0.001164461718427699

# Row 15: sCs -> FAPs (L1cam/Erbb3)
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "L1cam"
$ws.Range("C15").Value = "Erbb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.209654333333333
$ws.Range("H15").Value = 9.628962999999999
$ws.Range("I15").Value = 0.1209123959641148
$ws.Range("J15").Value = 0.1209123959641148
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.181585
$ws.Range("N15").Value = 0.544755
$ws.Range("O15").Value = 0.04185645645377991
$ws.Range("P15").Value = 0.04185645645377991
$ws.Range("Q15").Value = 0.5828250821183333
$ws.Range("R15").Value = 5.245425739064999
$ws.Range("S15").Value = 0.005060964436394165
$ws.Range("T15").Value = 0.005060964436394165

# Row 16: sCs -> M2 (L1cam/Erbb3)
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "L1cam"
$ws.Range("C16").Value = "Erbb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.209654333333333
$ws.Range("H16").Value = 9.628962999999999
$ws.Range("I16").Value = 0.1209123959641148
$ws.Range("J16").Value = 0.1209123959641148
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.4121513333333333
$ws.Range("N16").Value = 1.236454
$ws.Range("O16").Value = 0.0950034107224385
$ws.Range("P16").Value = 0.0950034107224385
$ws.Range("Q16").Value = 1.322863313022444
$ws.Range("R16").Value = 11.905769817202
$ws.Range("S16").Value = 0.01148709001521291
$ws.Range("T16").Value = 0.01148709001521291

# Row 17: sCs -> sCs (L1cam/Erbb3)
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "L1cam"
$ws.Range("C17").Value = "Erbb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 3.209654333333333
$ws.Range("H17").Value = 9.628962999999999
$ws.Range("I17").Value = 0.1209123959641148
$ws.Range("J17").Value = 0.1209123959641148
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.702762666666667
$ws.Range("N17").Value = 11.108288
$ws.Range("O17").Value = 0.8535095096842542
$ws.Range("P17").Value = 0.8535095096842543
$ws.Range("Q17").Value = 11.88458823837155
$ws.Range("R17").Value = 106.961294145344
$ws.Range("S17").Value = 0.10319987979408
$ws.Range("T17").Value = 0.10319987979408
